# break out stock.yaml completed
#
# 1) "day" sheet: D20 (bsecode for COLPAL) was captured as text ("500830")
#    instead of a number -> fix it to a real number.
# 2) "week" sheet: five new rows (47-51) were appended for the 18/06/2024
#    weekly run (ASTRAL, GLENMARK, BSOFT, GRANULES, BEL). The bsecode
#    column (D) for these rows was written as text, matching the source
#    feed for that day, so we force that column to Text format before
#    writing the values and then drop the style back to Normal so no
#    stray number-format lingers on the cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fix day!D20 bsecode to be numeric instead of text
# ---------------------------------------------------------------------
$dayWs = $wb.Worksheets.Item("day")
$dayWs.Cells.Item(20, 4).Value = 500830

# ---------------------------------------------------------------------
# 2) Append rows 47-51 to the "week" sheet
# ---------------------------------------------------------------------
$weekWs = $wb.Worksheets.Item("week")

$newRows = @(
    @{ Row=47; A=1; B="ASTRAL";    C="Astral Poly Technik Limited";      D="532830"; E=-0.54; F=2257.2;  G=360927;   H="week"; I="18/06/2024 11:32:29" },
    @{ Row=48; A=2; B="GLENMARK";  C="Glenmark Pharmaceuticals Limited"; D="532296"; E=0.46;  F=1242.85; G=1122110;  H="week"; I="18/06/2024 11:32:29" },
    @{ Row=49; A=3; B="BSOFT";     C="Birlasoft Ltd";                    D="532400"; E=0.86;  F=683.8;   G=3428629;  H="week"; I="18/06/2024 11:32:29" },
    @{ Row=50; A=4; B="GRANULES";  C="Granules India Limited";           D="532482"; E=2.65;  F=474;     G=1377669;  H="week"; I="18/06/2024 11:32:29" },
    @{ Row=51; A=5; B="BEL";       C="Bharat Electronics Limited";       D="500049"; E=2.79;  F=318.25;  G=64859245; H="week"; I="18/06/2024 11:32:29" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $weekWs.Cells.Item($row, 1).Value = $r.A
    $weekWs.Cells.Item($row, 2).Value = $r.B
    $weekWs.Cells.Item($row, 3).Value = $r.C

    # bsecode column: stored as text in the source feed even though it
    # looks numeric, so force Text format, write it, then reset the
    # style so the cell keeps no explicit style (matches rest of sheet).
    $dCell = $weekWs.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $r.D
    $dCell.Style = "Normal"

    $weekWs.Cells.Item($row, 5).Value = $r.E
    $weekWs.Cells.Item($row, 6).Value = $r.F
    $weekWs.Cells.Item($row, 7).Value = $r.G
    $weekWs.Cells.Item($row, 8).Value = $r.H
    $weekWs.Cells.Item($row, 9).Value = $r.I
}
